# Record a second data run (higher prominence) for the first data set, and
# add a couple of clarifying notes above the results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row just above the "distance:" definition (old row 17) for a
# new note about how prominence is now chosen/recorded per run.
$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = "Actually note, this way of choosing the prominence seemed not to work for first data run, so from now, will put a note about how prominence was selected for each run next to its row in the table. Probably going to try and aim for maybe the whole height span of visually looking roughly the widest noise part that doesn't seem to visually contain a peak."

# Insert another new row just above the blank line that precedes the results
# table (old row 21 after the first insert shifted it) for a second note.
$ws.Rows("21:21").Insert()
$ws.Range("A21").Value = "Also note, if things change between the runs, will try to note this next to data row."

# Append the new data run (second run of the first data set, with the
# prominence raised from 0.25 to 0.5) as a new row under the existing one.
$ws.Range("A25").Value = "sg_rr_20_025 2023-12-13 17-59-26.csv"
$ws.Range("B25").Value = 0.01
$ws.Range("C25").Value = 1000
$ws.Range("D25").Value = 5001
$ws.Range("E25").Value = 1530
$ws.Range("F25").Value = 1570
$ws.Range("G25").Value = 0.5
$ws.Range("H25").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 4.92
$ws.Range("K25").Value = 0.025354627641843101
$ws.Range("L25").Value = "Note I think Excel is not storing full number digits as printed out in Jupyter notebook, as I copy and paste the fsr mean and error into here, but it is keeping enough significant figures, and we should be able to look to GitHub or run data again."

# Vertically center the fsr_mean value of the new run.
$ws.Range("J25").VerticalAlignment = -4108

# Scroll the view down a bit and leave the selection just past the new row,
# matching where the author ended up after typing in the new data.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L26").Select()
